# Apply the "Sample Return Module Requirements" update:
#  - add a new worksheet "Sheet2" (after Sheet1) containing the SR module
#    requirements/constraints tables
#  - make the new sheet the active sheet
#  - move the selection on Sheet1 from B27 to B22

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Sheet1: update the remembered selection -----------------------------
$ws1.Range("B22").Select() | Out-Null

# ---- Add the new "Sheet2" worksheet (placed right after Sheet1) ----------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Column B is wide, like on Sheet1, to hold the long requirement text
$ws2.Columns.Item(2).ColumnWidth = 66.5

# ---- Fill in the new descriptive text, in authoring order, so that new
# ---- shared-string entries land at the same indices as the source file ---
$ws2.Range("B1").Value  = "Sample Return Module Requirements"
$ws2.Range("B6").Value  = "The Sample Return Arduino shall manage the operation of the motors"
$ws2.Range("B9").Value  = "The Sample Return Arduino shall manage the operation of the servos"
$ws2.Range("B16").Value = "The Sample Return Arduino will use a clockspring to connect to the primary module"
$ws2.Range("B15").Value = "The Sample Return Arduino will use the Sabertooth 2x25 for motor and servo control"
$ws2.Range("B17").Value = "The Sample Return Arduino will connect to the primary module through I2C"
$ws2.Range("B7").Value  = "The Arduino shall be capable of enabling and disabling the motors"
$ws2.Range("B4").Value  = "The Arduino shall parse packets properly"
$ws2.Range("B3").Value  = "The Arduino shall use two way data communications"
$ws2.Range("B2").Value  = "The Sample Return Arduino must be able to communicate with only the primary module"
$ws2.Range("A19").Value = "C5"
$ws2.Range("B5").Value  = "The Arduino will have a scheme to set the primary module to sample retrival mode"
$ws2.Range("B8").Value  = "The Arduino shall turn on or off a motor when a proper packet is sent"
$ws2.Range("B10").Value = "The Arduino shall be capable of enabling and disabling the servos"
$ws2.Range("B11").Value = "The Servo shall be programmed to run for a specified distance"
$ws2.Range("B19").Value = "The servos used will be part number: xxxxxxx"
$ws2.Range("B18").Value = "The motors used will be part number: xxxxxx"

# ---- Remaining "ID" / requirement-number cells (reuse existing strings) --
$ws2.Range("A1").Value  = "ID"
$ws2.Range("A2").Value  = "R1"
$ws2.Range("A3").Value  = "R1.1"
$ws2.Range("A4").Value  = "R1.2"
$ws2.Range("A5").Value  = "R1.3"
$ws2.Range("A6").Value  = "R2"
$ws2.Range("A7").Value  = "R2.1"
$ws2.Range("A8").Value  = "R2.2"
$ws2.Range("A9").Value  = "R3"
$ws2.Range("A10").Value = "R3.1"
$ws2.Range("A11").Value = "R3.2"

$ws2.Range("A14").Value = "ID"
$ws2.Range("B14").Value = "Constraints"
$ws2.Range("A15").Value = "C1"
$ws2.Range("A16").Value = "C2"
$ws2.Range("A17").Value = "C3"
$ws2.Range("A18").Value = "C4"

# ---- Formatting ------------------------------------------------------------
$ws2.Range("A1:B1").HorizontalAlignment = -4108
$ws2.Range("A1:B1").Font.Bold = $true
$ws2.Range("A14:B14").HorizontalAlignment = -4108
$ws2.Range("A14:B14").Font.Bold = $true

$ws2.Range("A3:A5").HorizontalAlignment = -4108
$ws2.Range("A7:A8").HorizontalAlignment = -4108
$ws2.Range("A10:A11").HorizontalAlignment = -4108
$ws2.Range("A12").HorizontalAlignment = -4108

# ---- Final view state: Sheet2 selected/active, cell A20 selected ---------
$ws2.Activate() | Out-Null
$ws2.Range("A20").Select() | Out-Null
